$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "Componente"

# 2. Insert three fresh rows right where the old "TLE4946-2L" row (10) used to live.
#    This shifts that whole row - including its distinctive thick-bottom-border row
#    formatting - down to row 13, where it will become the new final row
#    (RK09K1130AAU) once we overwrite its contents below. Rows 10-12 come back empty,
#    ready to receive the TLE4946-2L (restructured), LD-BZEN-1205 and DS18B20+ entries.
$ws.Rows("10:12").Insert(-4121) | Out-Null   # xlShiftDown

# 3. Quantities that were left blank before now get filled in.
$ws.Range("G8").Value = 12
$ws.Range("G9").Value = 5

# 4. Give the three new interior rows the same cell formatting as the row above
#    (row 9, BC857C.215) - border style, alignment, hyperlink font, etc.
$ws.Range("B9:G9").Copy() | Out-Null
$ws.Range("B10:G10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B11:G11").PasteSpecial(-4122) | Out-Null
$ws.Range("B12:G12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 5. Row 10 - TLE4946-2L: turn the old "-" placeholders into a real Mouser link and a
#    price.
$ws.Range("B10").Value = "TLE4946-2L"
$ws.Range("C10").Value = "Latch Hall de mare precizie"
$ws.Range("D10").Formula = '=HYPERLINK("https://ro.mouser.com/ProductDetail/Infineon-Technologies/TLE4946-2L?qs=sGAEpiMZZMvhQj7WZhFIAPDIuJuz6ELs2XiVbeCMsXM%3D","Mouser")'
$ws.Range("E10").Formula = '=HYPERLINK("https://www.infineon.com/dgdl/Infineon-TLE4946_2L-DS-v01_00-en.pdf?fileId=db3a304338ec6d390138fc8f905876d3","Link")'
$ws.Range("F10").Value = 3.76
$ws.Range("G10").Value = 1

# The price cell for this particular row lost its border in the original edit (it looks
# like it was typed directly over the old "-" cell without carrying the border along) -
# replicate that exactly.
$ws.Range("F10").Borders.LineStyle = 0
$ws.Range("F10").HorizontalAlignment = -4108   # xlCenter

# 6. Row 11 - new component: LD-BZEN-1205
$ws.Range("B11").Value = "LD-BZEN-1205"
$ws.Range("C11").Value = "Traductor de sunet: electromagnetic; fără generator încorporat"
$ws.Range("D11").Formula = '=HYPERLINK("https://www.tme.eu/ro/details/ld-bzen-1205/semnaliz-acust-electromag-fara-gener/loudity/","TME")'
$ws.Range("E11").Formula = '=HYPERLINK("https://www.tme.eu/Document/f21095d16ba59bdbd6f5bf5f2f188c94/ld-bzen-1205.pdf","Link")'
$ws.Range("F11").Value = 1.6033
$ws.Range("G11").Value = 1

# 7. Row 12 - new component: DS18B20+
$ws.Range("B12").Value = "DS18B20+"
$ws.Range("C12").Value = "Senzor temperatură; termometru digital; -55÷125°C; TO92; THT"
$ws.Range("D12").Formula = '=HYPERLINK("https://www.tme.eu/ro/details/ds18b20+/traductor-de-temperatura/maxim-integrated/","TME")'
$ws.Range("E12").Formula = '=HYPERLINK("https://www.tme.eu/Document/dd11228ca818d1ed1f3029ea47fef609/DS18B20+.pdf","Link")'
$ws.Range("F12").Value = 7.26
$ws.Range("G12").Value = 1

# 8. Row 13 (the row that got pushed down in step 2, still carrying the final/thick
#    border formatting) - new component: RK09K1130AAU
$ws.Range("B13").Value = "RK09K1130AAU"
$ws.Range("C13").Value = "Potenţiometru: axial; singură tură; 10kΩ; ±20%; 6mm; pt.PCB"
$ws.Range("D13").Formula = '=HYPERLINK("https://www.tme.eu/ro/details/rk09k1130aau/poten-axi-cu-pis-de-car-si-tura-uni/alps/","TME")'
$ws.Range("E13").Formula = '=HYPERLINK("https://www.tme.eu/Document/fb277d87df6cd625a5906cf03c82b793/ALPS_RK09K.PDF","Link")'
$ws.Range("F13").Value = 1.64
$ws.Range("G13").Value = 1

# D13 used to hold a plain string ("-") so it carries the non-hyperlink flavor of the
# row's style; now it holds a HYPERLINK formula like E13, so line its formatting up with
# E13's.
$ws.Range("E13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 9. Column D got a bit wider to accommodate the new "Mouser" link text.
$ws.Columns("D").ColumnWidth = 7.42578125

# 10. Match the saved cursor/selection position from the authored workbook.
$ws.Range("C23").Select() | Out-Null
